$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 135 and 136, shifting existing rows 135:160 down to 137:162
$ws.Range("A135:A136").EntireRow.Insert()

# --- New row 135 ---
$ws.Cells.Item(135, 1).Value = 6
$ws.Cells.Item(135, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(135, 3).Value = "Metropolitana"
$ws.Cells.Item(135, 4).Value = 44504
$ws.Cells.Item(135, 5).Value = 13
$ws.Cells.Item(135, 6).Value = 100112022
$ws.Cells.Item(135, 7).Value = "Arveja Verde"
$ws.Cells.Item(135, 8).Value = "Perfection"
$ws.Cells.Item(135, 9).Value = "Primera"
$ws.Cells.Item(135, 10).Value = 220
$ws.Cells.Item(135, 11).Value = 12000
$ws.Cells.Item(135, 12).Value = 13000
$ws.Cells.Item(135, 13).Value = 12545
$ws.Cells.Item(135, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(135, 15).Value = "Región Metropolitana"
$ws.Cells.Item(135, 16).Value = 502
$ws.Cells.Item(135, 17).Value = 25
$ws.Cells.Item(135, 18).Value = "Hortaliza"

# --- New row 136 ---
$ws.Cells.Item(136, 1).Value = 6
$ws.Cells.Item(136, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(136, 3).Value = "Metropolitana"
$ws.Cells.Item(136, 4).Value = 44504
$ws.Cells.Item(136, 5).Value = 13
$ws.Cells.Item(136, 6).Value = 100112022
$ws.Cells.Item(136, 7).Value = "Arveja Verde"
$ws.Cells.Item(136, 8).Value = "Sin especificar"
$ws.Cells.Item(136, 9).Value = "Primera"
$ws.Cells.Item(136, 10).Value = 350
$ws.Cells.Item(136, 11).Value = 10000
$ws.Cells.Item(136, 12).Value = 11000
$ws.Cells.Item(136, 13).Value = 10571
$ws.Cells.Item(136, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(136, 15).Value = "Región Metropolitana"
$ws.Cells.Item(136, 16).Value = 423
$ws.Cells.Item(136, 17).Value = 25
$ws.Cells.Item(136, 18).Value = "Hortaliza"
